# Applies the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Thu May  4 08:13:13 UTC 2023 with GitHub Actions"
#
# Column D ("Price") holds numeric-looking text (e.g. "326.62", "29.142.72")
# that must stay TEXT (the source cells are <is><t>, not <v> numbers). Typing
# such a value into a General-formatted cell would auto-convert it to a real
# number, so each Price cell is temporarily switched to Text format, written,
# then has that temporary format cleared again (cell keeps its text value/type,
# but does not end up with a stray number-format style applied).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    "D2" = '29.142.72'
    "D3" = '1.904.60'
    "D5" = '326.62'
    "D7" = '0.4628'
    "D8" = '0.3936'
    "D9" = '46.60'
    "D10" = '0.07952'
    "D11" = '0.9981'
    "D13" = '1.845.04'
    "D14" = '7.094'
    "D15" = '5.766'
    "D16" = '0.06945'
    "D17" = '88.43'
    "D18" = '1.003'
    "D20" = '17.13'
    "D21" = '1.002'
    "D22" = '29.168.55'
    "D23" = '5.358'
    "D24" = '11.10'
    "D25" = '2.092.67'
    "D26" = '2.064'
    "D27" = '156.54'
    "D29" = '6.011'
    "D31" = '119.00'
    "D32" = '0.09385'
    "D33" = '0.9197'
    "D34" = '5.333'
    "D35" = '1.349'
    "D36" = '3.259'
    "D37" = '1.203'
    "D38" = '0.05818'
    "D39" = '0.02101'
    "D40" = '7.954'
    "D41" = '0.5738'
    "D42" = '0.1796'
    "D43" = '9.958'
    "D44" = '11.98'
    "D45" = '0.5418'
    "D46" = '2.195'
    "D47" = '0.07088'
    "D49" = '2.552'
    "D50" = '112.03'
    "D51" = '1.042'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Coin name / link (row 37 <-> row 38 swap) and Volume(1h) percentages --
# these are already plain text (contain letters, "/", "%", spaces) so no
# special handling is required; Excel stores them as text automatically.
$otherUpdates = [ordered]@{
    "E2" = '  +1.48%  '
    "E3" = '  +1.86%  '
    "E4" = '  -0.03%  '
    "E5" = '  +0.69%  '
    "E6" = '  -0.04%  '
    "E7" = '  +0.51%  '
    "E8" = '  +2.03%  '
    "E9" = '  +0.80%  '
    "E10" = '  +1.37%  '
    "E11" = '  +1.75%  '
    "E12" = '  +2.14%  '
    "E13" = '  -2.51%  '
    "E14" = '  +1.50%  '
    "E15" = '  +1.25%  '
    "E16" = '  -0.14%  '
    "E17" = '  +0.16%  '
    "E18" = '  -0.10%  '
    "E19" = '  +0.39%  '
    "E20" = '  +2.33%  '
    "E21" = '  -0.05%  '
    "E22" = '  +1.55%  '
    "E23" = '  +1.73%  '
    "E24" = '  +0.26%  '
    "E25" = '  -1.16%  '
    "E26" = '  -1.61%  '
    "E27" = '  +2.80%  '
    "E28" = '  +1.25%  '
    "E29" = '  +2.28%  '
    "E30" = '  +0.26%  '
    "E31" = '  -0.06%  '
    "E32" = '  +0.68%  '
    "E33" = '  +0.33%  '
    "E34" = '  +0.76%  '
    "E35" = '  +1.13%  '
    "E36" = '  -1.88%  '
    "B37" = 'TrustWalletToken'
    "C37" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    "E37" = '  +5.23%  '
    "B38" = 'Hedera'
    "C38" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    "E38" = '  +0.72%  '
    "E39" = '  +1.32%  '
    "E40" = '  +4.01%  '
    "E42" = '  +0.73%  '
    "E43" = '  +2.16%  '
    "E44" = '  +1.77%  '
    "E45" = '  +2.45%  '
    "E46" = '  +3.81%  '
    "E47" = '  -1.75%  '
    "E48" = '  +2.26%  '
    "E49" = '  +5.95%  '
    "E50" = '  -0.84%  '
    "E51" = '  -7.17%  '
}

foreach ($addr in $otherUpdates.Keys) {
    $ws.Range($addr).Value = $otherUpdates[$addr]
}
